# Updates cryptos list values (Price / Volume(1h), and a few Coin/Link
# swaps) on Sheet1 to match the latest scrape, per commit:
# "Updated cryptos list on Sun Aug 25 02:41:38 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.158.80'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.00%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.758.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.20%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.98%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.59%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.31%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.609'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.35%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.112'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.96%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -13.68%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.389'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.48%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.158'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.71%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.246.93'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.20%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.15%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.089.48'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.04%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000154'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.49%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.760.34'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.58%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.21%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.65%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '361.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.19%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.81'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.55%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.559'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.95%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.81%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.78%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.172'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.14%  '
# Row 26
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.74%  '
# Row 27
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.32%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0939'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.45%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.19%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.31%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.01%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '168.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.39%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.42'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.52%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.96'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.99%  '
# Row 35
$ws.Range('E35').Value = '  +0.19%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.91%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.07%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.45%  '
# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.26'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.84%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.92%  '
# Row 41
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '333.57'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.11%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.30%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.96'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.00%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.97'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.22%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0596'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.51%  '
# Row 46
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0258'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.81%  '
# Row 47
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.639'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.98%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '136.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.60%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.102'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.81%  '
# Row 50
$ws.Range('E50').Value = '  +0.34%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.66%  '
